# Generate Report for Handoff
# Updates the localization-status workbook to mark c2e934a3 and d4ffe16f
# files as "Ready for handoff" and refresh handoff timestamps / error detail.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: rows 4 (c2e934a3) and 5 (d4ffe16f) ---
$overview.Range("E4").Value = "Ready for handoff"
$overview.Range("F4").Value = "Ready for handoff"
$overview.Range("G4").Value = "2016-08-28 06:31:21"

$overview.Range("E5").Value = "Ready for handoff"
$overview.Range("F5").Value = "Ready for handoff"
$overview.Range("G5").Value = "2016-08-28 06:31:21"

# --- zh-cn sheet: rows 4 (c2e934a3) and 5 (d4ffe16f) ---
$zhcn.Range("C4").Value = "Ready for handoff"
$zhcn.Range("H4").Value = "2016-08-28 06:31:17"
$zhcn.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a628533bfac68649978fe8346b8dbc5237dc88f3/e2e/c2e934a3-4f62-486d-96fa-7abd9d764473.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c202cd3c793a03ff422d3fc305dfbbad9855f833/e2e/c2e934a3-4f62-486d-96fa-7abd9d764473.md."

$zhcn.Range("C5").Value = "Ready for handoff"
$zhcn.Range("H5").Value = "2016-08-28 06:31:17"
$zhcn.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a628533bfac68649978fe8346b8dbc5237dc88f3/e2e/d4ffe16f-491c-4f8b-8185-92b9e07263f5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c202cd3c793a03ff422d3fc305dfbbad9855f833/e2e/d4ffe16f-491c-4f8b-8185-92b9e07263f5.md."

$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667

# --- de-de sheet: rows 4 (c2e934a3) and 5 (d4ffe16f) ---
$dede.Range("C4").Value = "Ready for handoff"
$dede.Range("H4").Value = "2016-08-28 06:31:21"
$dede.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a628533bfac68649978fe8346b8dbc5237dc88f3/e2e/c2e934a3-4f62-486d-96fa-7abd9d764473.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c202cd3c793a03ff422d3fc305dfbbad9855f833/e2e/c2e934a3-4f62-486d-96fa-7abd9d764473.md."

$dede.Range("C5").Value = "Ready for handoff"
$dede.Range("H5").Value = "2016-08-28 06:31:21"
$dede.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a628533bfac68649978fe8346b8dbc5237dc88f3/e2e/d4ffe16f-491c-4f8b-8185-92b9e07263f5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c202cd3c793a03ff422d3fc305dfbbad9855f833/e2e/d4ffe16f-491c-4f8b-8185-92b9e07263f5.md."

$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
